$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1000.63635
$ws.Range("I15").Value = 1000.63635
$ws.Range("K15").Value = 3001.90905
$ws.Range("M15").Value = -2832.90905
$ws.Range("H125").Value = 2850
$ws.Range("I125").Value = 2850
$ws.Range("K125").Value = 25650
$ws.Range("M125").Value = -23190
$ws.Range("H137").Value = 5832.7144
$ws.Range("I137").Value = 5832.7144
$ws.Range("K137").Value = 17498.1432
$ws.Range("M137").Value = -14948.1432
$ws.Range("H138").Value = 4221.273
$ws.Range("J138").Value = 4937.1113
$ws.Range("L138").Value = 14811.3339
$ws.Range("N138").Value = -25091.3339
$ws.Range("H141").Value = 1166.6666
$ws.Range("I141").Value = 1250
$ws.Range("K141").Value = 3750
$ws.Range("M141").Value = 1430
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5938.2
$ws.Range("I74").Value = 5938.2
$ws.Range("K74").Value = 5938.2
$ws.Range("M74").Value = -5064.2
$ws.Range("H77").Value = 5938.2
$ws.Range("I77").Value = 5938.2
$ws.Range("K77").Value = 29691
$ws.Range("M77").Value = -25323
$ws.Range("H132").Value = 2474.6
$ws.Range("I132").Value = 1056.9333
$ws.Range("K132").Value = 3170.7999
$ws.Range("M132").Value = -640.7999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2558.5293
$ws.Range("J134").Value = 3949.5
$ws.Range("L134").Value = 11848.5
$ws.Range("N134").Value = -16918.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3000.8
$ws.Range("I31").Value = 1726.2858
$ws.Range("J31").Value = 3687.077
$ws.Range("K31").Value = 1726.2858
$ws.Range("L31").Value = 3687.077
$ws.Range("M31").Value = -1431.2858
$ws.Range("N31").Value = -4277.077
$ws.Range("H34").Value = 3000.8
$ws.Range("I34").Value = 1726.2858
$ws.Range("J34").Value = 3687.077
$ws.Range("K34").Value = 1726.2858
$ws.Range("L34").Value = 3687.077
$ws.Range("M34").Value = -1524.2858
$ws.Range("N34").Value = -4091.077
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 74000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 74000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 74000
$ws.Range("N130").Value = -84040
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 5259.385
$ws.Range("I132").Value = 5137.4
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 15412.2
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -12882.2
$ws.Range("N132").Value = -22058
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 5299
$ws.Range("I134").Value = 5942.706
$ws.Range("J134").Value = 1651.3334
$ws.Range("K134").Value = 17828.118
$ws.Range("L134").Value = 4954.0002
$ws.Range("M134").Value = -15293.118
$ws.Range("N134").Value = -10024.0002
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 75000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 75000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 213.85715
$ws.Range("I86").Value = 209.25
$ws.Range("J86").Value = 220
$ws.Range("K86").Value = 627.75
$ws.Range("L86").Value = 660
$ws.Range("M86").Value = 558.25
$ws.Range("N86").Value = -3032
$ws.Range("H89").Value = 213.85715
$ws.Range("I89").Value = 209.25
$ws.Range("J89").Value = 220
$ws.Range("K89").Value = 1883.25
$ws.Range("L89").Value = 1980
$ws.Range("M89").Value = 4044.75
$ws.Range("N89").Value = -13836
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 1148.75
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1148.75
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3446.25
$ws.Range("N121").Value = -6066.25
$ws.Range("H122").Value = 933.1111
$ws.Range("I122").Value = 366.66666
$ws.Range("J122").Value = 2066
$ws.Range("K122").Value = 3299.99994
$ws.Range("L122").Value = 18594
$ws.Range("M122").Value = -849.9999399999997
$ws.Range("N122").Value = -23494
$ws.Range("H123").Value = 4666.5557
$ws.Range("I123").Value = 1999
$ws.Range("J123").Value = 5000
$ws.Range("K123").Value = 5997
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = -3547
$ws.Range("N123").Value = -19900
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 6249.75
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 6666.3335
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 19999.0005
$ws.Range("M126").Value = -10060
$ws.Range("N126").Value = -29879.0005
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 240000
$ws.Range("I128").Value = 240000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 720000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -715020
$ws.Range("H129").Value = 1433.3334
$ws.Range("I129").Value = 1433.3334
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 4300.0002
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 699.9997999999996
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 2214.5
$ws.Range("I131").Value = 1727.25
$ws.Range("J131").Value = 2701.75
$ws.Range("K131").Value = 5181.75
$ws.Range("L131").Value = 8105.25
$ws.Range("M131").Value = -141.75
$ws.Range("N131").Value = -18185.25
$ws.Range("H132").Value = 194
$ws.Range("I132").Value = 194
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1746
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 784
$ws.Range("H133").Value = 3998
$ws.Range("I133").Value = 3998
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 11994
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -6934
$ws.Range("H134").Value = 5332.5
$ws.Range("I134").Value = 4399
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 13197
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -8127
$ws.Range("N134").Value = -40140
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 4500
$ws.Range("N137").Value = -14700
$ws.Range("H138").Value = 4653
$ws.Range("I138").Value = 1980
$ws.Range("J138").Value = 5989.5
$ws.Range("K138").Value = 5940
$ws.Range("L138").Value = 17968.5
$ws.Range("M138").Value = -800
$ws.Range("N138").Value = -28248.5
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 4994.1113
$ws.Range("I140").Value = 5393.375
$ws.Range("J140").Value = 1800
$ws.Range("K140").Value = 16180.125
$ws.Range("L140").Value = 5400
$ws.Range("M140").Value = -11000.125
$ws.Range("N140").Value = -15760
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2180
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1145.091
$ws.Range("I136").Value = 1109.6
$ws.Range("K136").Value = 3328.8
$ws.Range("M136").Value = -778.7999999999997
